# Add UMI fields to Excel report
#
# The "schemas + fields" sheet lists field names alphabetically down column A,
# with a checkmark in the column for every assay/schema that uses that field
# (header row 1: B=af, C=antibodies, ... V=scrnaseq, ...).
#
# Three new UMI-related fields belong (alphabetically) right before
# "uniprot_accession_number" (currently row 222), and they apply to the
# "scrnaseq" schema (column V):
#   umi_offset
#   umi_read
#   umi_size

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 222, pushing the existing rows (uniprot_accession_number
# onward) down to make room, without disturbing any other formatting/content.
$ws.Range("A222:A224").EntireRow.Insert()

$checkmark = [char]0x2713

$ws.Range("A222").Value = "umi_offset"
$ws.Range("V222").Value = $checkmark

$ws.Range("A223").Value = "umi_read"
$ws.Range("V223").Value = $checkmark

$ws.Range("A224").Value = "umi_size"
$ws.Range("V224").Value = $checkmark
